$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 150
$ws.Range("J2").Value = 150
$ws.Range("L2").Value = 150
$ws.Range("N2").Value = -376
$ws.Range("H6").Value = 8865.700000000001
$ws.Range("I6").Value = 8865.700000000001
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 26597.1
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -26485.1
$ws.Range("N6").Value = ""
$ws.Range("H16").Value = 3067.5557
$ws.Range("I16").Value = 2921.6
$ws.Range("K16").Value = 2921.6
$ws.Range("M16").Value = -2691.6
$ws.Range("H40").Value = 4273.857
$ws.Range("I40").Value = 2224.5
$ws.Range("K40").Value = 2224.5
$ws.Range("M40").Value = -2049.5
$ws.Range("H43").Value = 5312
$ws.Range("I43").Value = 2135
$ws.Range("J43").Value = 9759.799999999999
$ws.Range("K43").Value = 2135
$ws.Range("L43").Value = 9759.799999999999
$ws.Range("M43").Value = -2066
$ws.Range("N43").Value = -9897.799999999999
$ws.Range("H53").Value = 381.57144
$ws.Range("I53").Value = 224.3077
$ws.Range("K53").Value = 224.3077
$ws.Range("M53").Value = 412.6923
$ws.Range("H98").Value = 47627244
$ws.Range("I98").Value = 62506670
$ws.Range("K98").Value = 62506670
$ws.Range("M98").Value = -62505172
$ws.Range("H100").Value = 5509.4443
$ws.Range("I100").Value = 4798.143
$ws.Range("J100").Value = 7999
$ws.Range("K100").Value = 4798.143
$ws.Range("L100").Value = 7999
$ws.Range("M100").Value = -4257.143
$ws.Range("N100").Value = -9081
$ws.Range("H112").Value = 5200
$ws.Range("J112").Value = 5250
$ws.Range("L112").Value = 15750
$ws.Range("N112").Value = -17966
$ws.Range("H113").Value = 18218.75
$ws.Range("I113").Value = 4291.6665
$ws.Range("J113").Value = 60000
$ws.Range("K113").Value = 4291.6665
$ws.Range("L113").Value = 60000
$ws.Range("M113").Value = -1037.6665
$ws.Range("N113").Value = -66508
$ws.Range("H122").Value = 47627244
$ws.Range("I122").Value = 62506670
$ws.Range("K122").Value = 187520010
$ws.Range("M122").Value = -187517560
$ws.Range("H132").Value = 2399.0667
$ws.Range("I132").Value = 2106.1428
$ws.Range("K132").Value = 6318.428400000001
$ws.Range("M132").Value = -3788.428400000001
$ws.Range("H137").Value = 6457
$ws.Range("I137").Value = 3558.8
$ws.Range("J137").Value = 8268.375
$ws.Range("K137").Value = 10676.4
$ws.Range("L137").Value = 24805.125
$ws.Range("M137").Value = -8126.400000000001
$ws.Range("N137").Value = -29905.125
$ws.Range("H138").Value = 2522.738
$ws.Range("J138").Value = 3024.375
$ws.Range("L138").Value = 9073.125
$ws.Range("N138").Value = -19353.125

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 58138
$ws.Range("I31").Value = 4999.6665
$ws.Range("J31").Value = 111276.336
$ws.Range("K31").Value = 4999.6665
$ws.Range("L31").Value = 111276.336
$ws.Range("M31").Value = -4705.6665
$ws.Range("N31").Value = -111864.336
$ws.Range("H32").Value = 11910725
$ws.Range("I32").Value = 12826847
$ws.Range("J32").Value = 1139.6666
$ws.Range("K32").Value = 12826847
$ws.Range("L32").Value = 1139.6666
$ws.Range("M32").Value = -12826560
$ws.Range("N32").Value = -1713.6666
$ws.Range("H74").Value = 11914503
$ws.Range("I74").Value = 22728838
$ws.Range("J74").Value = 18733.1
$ws.Range("K74").Value = 22728838
$ws.Range("L74").Value = 18733.1
$ws.Range("M74").Value = -22727964
$ws.Range("N74").Value = -20481.1
$ws.Range("H77").Value = 11914503
$ws.Range("I77").Value = 22728838
$ws.Range("J77").Value = 18733.1
$ws.Range("K77").Value = 113644190
$ws.Range("L77").Value = 93665.5
$ws.Range("M77").Value = -113639822
$ws.Range("N77").Value = -102401.5
$ws.Range("H102").Value = 65825.664
$ws.Range("I102").Value = 65825.664
$ws.Range("K102").Value = 65825.664
$ws.Range("M102").Value = -64203.664
$ws.Range("H132").Value = 5297
$ws.Range("I132").Value = 2200.5122
$ws.Range("K132").Value = 6601.5366
$ws.Range("M132").Value = -4071.5366

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2367.5881
$ws.Range("I99").Value = 1987.3334
$ws.Range("J99").Value = 3280.2
$ws.Range("K99").Value = 1987.3334
$ws.Range("L99").Value = 3280.2
$ws.Range("M99").Value = -489.3334
$ws.Range("N99").Value = -6276.2

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 557677.9399999999
$ws.Range("I31").Value = 1802.4445
$ws.Range("J31").Value = 974584.5600000001
$ws.Range("K31").Value = 1802.4445
$ws.Range("L31").Value = 974584.5600000001
$ws.Range("M31").Value = -1507.4445
$ws.Range("N31").Value = -975174.5600000001
$ws.Range("H34").Value = 557677.9399999999
$ws.Range("I34").Value = 1802.4445
$ws.Range("J34").Value = 974584.5600000001
$ws.Range("K34").Value = 1802.4445
$ws.Range("L34").Value = 974584.5600000001
$ws.Range("M34").Value = -1600.4445
$ws.Range("N34").Value = -974988.5600000001

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 999
$ws.Range("I23").Value = 998.5
$ws.Range("K23").Value = 2995.5
$ws.Range("M23").Value = -2760.5
$ws.Range("I51").Value = 5570.5713
$ws.Range("J51").Value = 22102.5
$ws.Range("K51").Value = 16711.7139
$ws.Range("L51").Value = 66307.5
$ws.Range("M51").Value = -16251.7139
$ws.Range("N51").Value = -67227.5
$ws.Range("H68").Value = 2896
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 2896
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4018.5
$ws.Range("I102").Value = 4387.273
$ws.Range("J102").Value = 2666.3333
$ws.Range("K102").Value = 4387.273
$ws.Range("L102").Value = 2666.3333
$ws.Range("M102").Value = -2765.273
$ws.Range("N102").Value = -5910.3333
$ws.Range("H132").Value = 111114510
$ws.Range("I132").Value = 142860260
$ws.Range("K132").Value = 428580780
$ws.Range("M132").Value = -428578250

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 111113310
$ws.Range("I93").Value = 111113310
$ws.Range("K93").Value = 111113310
$ws.Range("M93").Value = -111112062
$ws.Range("H102").Value = 69755.60000000001
$ws.Range("I102").Value = 69989
$ws.Range("J102").Value = 69697.25
$ws.Range("K102").Value = 69989
$ws.Range("L102").Value = 69697.25
$ws.Range("M102").Value = -66744
$ws.Range("N102").Value = -76187.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7666
$ws.Range("I96").Value = 1000
$ws.Range("K96").Value = 1000
$ws.Range("M96").Value = 373
$ws.Range("H122").Value = 5551.625
$ws.Range("I122").Value = 3727.3845
$ws.Range("J122").Value = 7707.5454
$ws.Range("K122").Value = 11182.1535
$ws.Range("L122").Value = 23122.6362
$ws.Range("M122").Value = -8732.1535
$ws.Range("N122").Value = -28022.6362

